$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook tracks localization status for a set of files. Two files,
#   93e4a31f-1314-42f3-8e56-bb922eefd6e8  (currently row 4 on every sheet)
#   99ead3f0-7917-49df-9654-a4f94c45a773  (currently row 5 on every sheet)
# need to trade places: 99ead3f0 moves up into row 4 (and its status flips
# from "Ready for handoff" to "In Translation", reflecting that the file was
# sent back for further translation), while 93e4a31f moves down into row 5
# keeping its "Ready for handoff" status.
#
# Across all three worksheets (Overview, zh-cn, de-de) we:
#   1. swap the display text of the existing hyperlinks that live on rows 4/5
#      (without touching which URL they point at - that stays fixed to the
#      worksheet row/column position), and
#   2. swap/update the actual cell values so the visible text matches.
# ---------------------------------------------------------------------------

# --- Sheet 1: "Overview" ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Address -like "*93e4a31f-1314-42f3-8e56-bb922eefd6e8.md") {
        $h.TextToDisplay = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
    } elseif ($h.Address -like "*99ead3f0-7917-49df-9654-a4f94c45a773.md") {
        $h.TextToDisplay = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
    }
}

$wsOverview.Range("A4").Value = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

$wsOverview.Range("A5").Value = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

# --- Sheet 2: "zh-cn" --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Address -like "*93e4a31f-1314-42f3-8e56-bb922eefd6e8.md") {
        $h.TextToDisplay = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
    } elseif ($h.Address -like "*99ead3f0-7917-49df-9654-a4f94c45a773.md") {
        $h.TextToDisplay = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
    } elseif ($h.Address -like "*93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.zh-cn.xlf") {
        $h.TextToDisplay = "99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.zh-cn.xlf"
    } elseif ($h.Address -like "*99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.zh-cn.xlf") {
        $h.TextToDisplay = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.zh-cn.xlf"
    }
}

$wsZhCn.Range("A4").Value = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
$wsZhCn.Range("B4").Value = "In Translation"
$wsZhCn.Range("C4").Value = "99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2016-02-25 06:45:03"

$wsZhCn.Range("A5").Value = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
$wsZhCn.Range("B5").Value = "Ready for handoff"
$wsZhCn.Range("C5").Value = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.zh-cn.xlf"
$wsZhCn.Range("D5").Value = "2016-02-25 06:45:48"

# --- Sheet 3: "de-de" --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Address -like "*93e4a31f-1314-42f3-8e56-bb922eefd6e8.md") {
        $h.TextToDisplay = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
    } elseif ($h.Address -like "*99ead3f0-7917-49df-9654-a4f94c45a773.md") {
        $h.TextToDisplay = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
    } elseif ($h.Address -like "*93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.de-de.xlf") {
        $h.TextToDisplay = "99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.de-de.xlf"
    } elseif ($h.Address -like "*99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.de-de.xlf") {
        $h.TextToDisplay = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.de-de.xlf"
    }
}

$wsDeDe.Range("A4").Value = "99ead3f0-7917-49df-9654-a4f94c45a773.md"
$wsDeDe.Range("B4").Value = "In Translation"
$wsDeDe.Range("C4").Value = "99ead3f0-7917-49df-9654-a4f94c45a773.7b5ccef2be4187d8d14e2b451a398be9125d803f.de-de.xlf"
$wsDeDe.Range("D4").Value = "2016-02-25 06:45:17"

$wsDeDe.Range("A5").Value = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.md"
$wsDeDe.Range("B5").Value = "Ready for handoff"
$wsDeDe.Range("C5").Value = "93e4a31f-1314-42f3-8e56-bb922eefd6e8.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.de-de.xlf"
$wsDeDe.Range("D5").Value = "2016-02-25 06:46:01"
